$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting existing rows 56-74 down to 57-75.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record (values copied from
# the same static columns used throughout this block, dates/values per diff).
$ws.Cells.Item(56, 1).Value = 5
$ws.Cells.Item(56, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(56, 3).Value = "Maule"
$ws.Cells.Item(56, 4).Value = 44795
$ws.Cells.Item(56, 5).Value = 7
$ws.Cells.Item(56, 6).Value = 100112026
$ws.Cells.Item(56, 7).Value = "Haba"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 150
$ws.Cells.Item(56, 11).Value = 10000
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 10000
$ws.Cells.Item(56, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Región del Maule"
$ws.Cells.Item(56, 16).Value = 400
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"
